# Applies the "Added new plotting data from ADS" commit:
#  - Adds a new "S11_ADS_FEM" series (row + column) to the correlation
#    matrices on MLIN_Disc_Monopole and CPW_Disc_Monopole.
#  - Updates the S11_ADS correlation values on Sierp_Patch_3rd and
#    Sierp_Patch_2rd (recomputed correlation matrices).
#  - Updates the active selection on Planilha1.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- MLIN_Disc_Monopole: add S11_ADS_FEM row/column ---------------------
$ws10 = $wb.Worksheets.Item("MLIN_Disc_Monopole")

$ws10.Range("E1").Value2 = "S11_ADS_FEM"
$ws10.Range("D1").Copy()
$ws10.Range("E1").PasteSpecial($xlPasteFormats)

$ws10.Range("A5").Value2 = "S11_ADS_FEM"
$ws10.Range("A4").Copy()
$ws10.Range("A5").PasteSpecial($xlPasteFormats)

$ws10.Range("E2").Value2 = -0.6501231965313155
$ws10.Range("E3").Value2 = 1
$ws10.Range("B5").Value2 = -0.6501231965313155
$ws10.Range("C5").Value2 = 1
$ws10.Range("E5").Value2 = 1

$excel.CutCopyMode = 0

# --- CPW_Disc_Monopole: add S11_ADS_FEM row/column -----------------------
$ws5 = $wb.Worksheets.Item("CPW_Disc_Monopole")

$ws5.Range("E1").Value2 = "S11_ADS_FEM"
$ws5.Range("D1").Copy()
$ws5.Range("E1").PasteSpecial($xlPasteFormats)

$ws5.Range("A5").Value2 = "S11_ADS_FEM"
$ws5.Range("A4").Copy()
$ws5.Range("A5").PasteSpecial($xlPasteFormats)

$ws5.Range("E2").Value2 = -0.4036618505910266
$ws5.Range("E3").Value2 = 0.9999999999999999
$ws5.Range("B5").Value2 = -0.4036618505910266
$ws5.Range("C5").Value2 = 0.9999999999999999
$ws5.Range("E5").Value2 = 1

$excel.CutCopyMode = 0

# --- Sierp_Patch_3rd: recomputed correlation values ----------------------
$ws6 = $wb.Worksheets.Item("Sierp_Patch_3rd")
$ws6.Range("C2").Value2 = -0.5015996501889113
$ws6.Range("B3").Value2 = -0.5015996501889113
$ws6.Range("D3").Value2 = -0.3896835492505257
$ws6.Range("C4").Value2 = -0.3896835492505257

# --- Sierp_Patch_2rd: recomputed correlation values ----------------------
$ws8 = $wb.Worksheets.Item("Sierp_Patch_2rd")
$ws8.Range("C2").Value2 = 0.4479707490939402
$ws8.Range("B3").Value2 = 0.4479707490939402
$ws8.Range("D3").Value2 = 0.8673387799484744
$ws8.Range("C4").Value2 = 0.8673387799484744

# --- Planilha1: update active selection -----------------------------------
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws1.Activate()
$ws1.Range("G25").Select()
